# Generate Report for Handback
# Removes the stale "d65a6fcc-...md" handback row (row 3) from every sheet
# and refreshes the "Correspond Handoff/Handback Datetime" timestamps on the
# remaining row for the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview": drop row 3 (the d65a6fcc... file) entirely.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Rows.Item(3).Delete()

# Row-3 deletion leaves a dangling hyperlink pointing at the now-removed A3;
# rebuild the hyperlink collection so only the surviving row's link remains.
$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/28ae6ea4a2260ee1f8454732c348f02e95df413c/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn": drop row 3, bump the handoff/handback datetimes on row 2,
# then rebuild the surviving hyperlinks.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Rows.Item(3).Delete()

$zhcn.Range("E2").Value = "2016-03-21 04:45:43"
$zhcn.Range("H2").Value = "2016-03-21 04:46:04"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/28ae6ea4a2260ee1f8454732c348f02e95df413c/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.md")
$zhcn.Hyperlinks.Add($zhcn.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/28ae6ea4a2260ee1f8454732c348f02e95df413c/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", ".md")
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35c5577c1421a04dd00943ffe4f9ff46b46d20ad/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fdd845ca30516e1ec02dd0295f021736166f8373/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bc7b6df3be8c920f5ca808e426cebb9fa6d68c47/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de": same treatment as zh-cn, with the de-de specific timestamps
# and hyperlink targets.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Rows.Item(3).Delete()

$dede.Range("E2").Value = "2016-03-21 04:45:46"
$dede.Range("H2").Value = "2016-03-21 04:46:10"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/28ae6ea4a2260ee1f8454732c348f02e95df413c/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.md")
$dede.Hyperlinks.Add($dede.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/28ae6ea4a2260ee1f8454732c348f02e95df413c/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", ".md")
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6842789a229f253871ea1b5b648aa7902ccc532b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4114aca3226019f2e1b8aaa58097cfbe99688b4f/e2e/4f079f8c-42cb-47c3-ad76-af83946074ac.md", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/6ebf4dd522bfdf71a8a0020248f8bfa486f80a1b/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf", "", "", "4f079f8c-42cb-47c3-ad76-af83946074ac.ee1e8986c3b4725fba3fae01e0aafc1fbfcb907d.de-de.xlf")
